# ORT dummy scan report: rename the scanner referenced in the issue
# messages from "FileCounter" to "Dummy" (the FileCounter scanner was
# removed in favor of a Dummy scanner used by the new
# ScannerIntegrationFunTest), and restore the view/selection state that
# Excel records after making that edit interactively.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Summary")
$ws2 = $wb.Worksheets.Item("Gradle org.ossreviewtoolkit.gra")

# --- Content edit -----------------------------------------------------
# Update the "Gradle org.ossreviewtoolkit.gra" sheet first, then
# "Summary" -- this matches the order the shared strings end up stored
# in when the workbook is saved.
$old2 = $ws2.Range("F11").Text
$ws2.Range("F11").Value = $old2.Replace("FileCounter", "Dummy")
$ws2.Rows.Item(11).RowHeight = 15

$old1 = $ws1.Range("F11").Text
$ws1.Range("F11").Value = $old1.Replace("FileCounter", "Dummy")
$ws1.Rows.Item(11).RowHeight = 30

# --- View / selection state --------------------------------------------
# The Gradle sheet's cursor ends on F11, Summary's cursor ends on F20,
# and Summary is left as the active (selected) tab.
$ws2.Activate()
$ws2.Range("F11").Select()

$ws1.Activate()
$ws1.Range("F20").Select()

Write-Host "Replaced FileCounter with Dummy in both sheets."
